$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.356.91'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.934.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.11%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.03'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7103'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9992'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3272'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -9.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.51'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06863'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8059'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08085'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.933.63'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.425'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.41'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.54'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '259.78'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.348.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.58%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008030'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.809'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.187.74'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.898'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.728'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.94'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.361'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.10'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.51%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1331'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.26%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.560'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.350'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.436'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.228'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05093'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.228'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7454'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.735'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01974'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.831'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.79'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.582'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4478'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.002'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -9.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8367'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.02'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.829'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.309'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.46'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.486'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4100'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.96%  '
